# Auto-applies the per-cell text updates from the commit diff.
# Column D ("Price") cells that look like plain decimal numbers get a
# leading apostrophe so Excel keeps them as literal text (matching the
# workbook's existing inlineStr/text storage) instead of coercing them
# to numbers and silently dropping significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.384.87"
$ws.Range("E2").Value = "  -2.59%  "
$ws.Range("D3").Value = "3.452.69"
$ws.Range("E3").Value = "  +3.36%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'236.76"
$ws.Range("E5").Value = "  -5.46%  "
$ws.Range("D6").Value = "'637.13"
$ws.Range("E6").Value = "  -2.91%  "
$ws.Range("D7").Value = "'1.43"
$ws.Range("E7").Value = "  +1.44%  "
$ws.Range("D8").Value = "'0.397"
$ws.Range("E8").Value = "  -6.35%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'0.972"
$ws.Range("E10").Value = "  -3.77%  "
$ws.Range("D11").Value = "3.449.74"
$ws.Range("E11").Value = "  +3.37%  "
$ws.Range("D12").Value = "'42.03"
$ws.Range("E12").Value = "  +3.11%  "
$ws.Range("E13").Value = "  -4.70%  "
$ws.Range("D14").Value = "'6.17"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "94.324.15"
$ws.Range("E15").Value = "  -2.40%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "4.101.32"
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("D17").Value = "'0.0000252"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "'8.36"
$ws.Range("E18").Value = "  -4.12%  "
$ws.Range("D19").Value = "3.447.43"
$ws.Range("E19").Value = "  +3.20%  "
$ws.Range("D20").Value = "'17.62"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").Value = "'11.34"
$ws.Range("E21").Value = "  +6.03%  "
$ws.Range("D22").Value = "'0.500"
$ws.Range("E22").Value = "  -11.62%  "
$ws.Range("D23").Value = "'498.65"
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("D24").Value = "'3.13"
$ws.Range("E24").Value = "  -6.26%  "
$ws.Range("D25").Value = "'6.60"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").Value = "'0.0000189"
$ws.Range("E26").Value = "  -4.84%  "
$ws.Range("D27").Value = "'91.23"
$ws.Range("E27").Value = "  -5.61%  "
$ws.Range("D28").Value = "3.639.53"
$ws.Range("E28").Value = "  +3.47%  "
$ws.Range("D29").Value = "'11.93"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("D30").Value = "'11.71"
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("D31").Value = "'0.993"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").Value = "'2.73"
$ws.Range("E32").Value = "  +7.61%  "
$ws.Range("D33").Value = "'0.135"
$ws.Range("E33").Value = "  -7.47%  "
$ws.Range("D34").Value = "'0.182"
$ws.Range("E34").Value = "  -3.37%  "
$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "'30.25"
$ws.Range("E36").Value = "  +6.46%  "
$ws.Range("D37").Value = "'0.563"
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("D38").Value = "'535.44"
$ws.Range("E38").Value = "  +5.54%  "
$ws.Range("D39").Value = "'7.62"
$ws.Range("E39").Value = "  -2.67%  "
$ws.Range("D40").Value = "'1.44"
$ws.Range("E40").Value = "  -4.44%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.151"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'0.931"
$ws.Range("E42").Value = "  +11.21%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").Value = "'1.68"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").Value = "'5.55"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").Value = "'0.0411"
$ws.Range("E47").Value = "  -5.38%  "
$ws.Range("B48").Value = "MantraDAO"
$ws.Range("C48").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D48").Value = "'3.50"
$ws.Range("E48").Value = "  -5.13%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'2.15"
$ws.Range("E49").Value = "  +8.13%  "
$ws.Range("D50").Value = "'53.47"
$ws.Range("E50").Value = "  -2.53%  "
$ws.Range("D51").Value = "'3.18"
$ws.Range("E51").Value = "  +2.26%  "
